$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 275; this shifts the existing rows 275-296
# down to 276-297 (old row 296 becomes new row 297) and grows the used
# range from A1:R296 to A1:R297.
$ws.Rows(275).Insert()

# Populate the newly inserted row 275 with the new weekly data point.
$ws.Range("A275").Value = 3
$ws.Range("B275").Value = "Femacal de La Calera"
$ws.Range("C275").Value = "Coquimbo"
$ws.Range("D275").Value = 44714
$ws.Range("E275").Value = 5
$ws.Range("F275").Value = 100112001
$ws.Range("G275").Value = "Berenjena"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 150
$ws.Range("K275").Value = 6500
$ws.Range("L275").Value = 7000
$ws.Range("M275").Value = 6733
$ws.Range("N275").Value = "$/caja 60 unidades"
$ws.Range("O275").Value = "Región de Arica y Parinacota"
$ws.Range("P275").Value = 112
$ws.Range("Q275").Value = 60
$ws.Range("R275").Value = "Hortaliza"
